$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# SignIn sheet: selection moves back to A1 (was A3)
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SignIn")
$ws1.Range("A1").Select()

# -----------------------------------------------------------------
# ShareSkill sheet: edit A2, add row 3, add Category/Subcategory cols
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ShareSkill")

$ws2.Range("A2").Value = "Title is Entered"

$ws2.Range("A3").Value = "Title editted"
$ws2.Range("B3").Value = "Description editted"

# New header cells O1/P1 - copy header fill style from an existing header cell
$ws2.Range("D1").Copy($ws2.Range("O1"))
$ws2.Range("O1").Value = "Category"
$ws2.Range("D1").Copy($ws2.Range("P1"))
$ws2.Range("P1").Value = "Subcategory"

$ws2.Range("O2").Value = "Graphics & Design"
$ws2.Range("P2").Value = "Book & Album covers"

$ws2.Columns.Item(15).ColumnWidth = 17.6
$ws2.Columns.Item(16).ColumnWidth = 30

$ws2.Activate()
$ws2.Range("P10").Select()

# -----------------------------------------------------------------
# New sheet: ManageListing
# -----------------------------------------------------------------
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $wsLast)
$ws3.Name = "ManageListing"

$ws3.Range("A1").Value = "url"
$ws3.Range("B1").Value = "x"
$ws3.Range("C1").Value = "y"
$ws2.Range("D1").Copy($ws3.Range("B1"))
$ws3.Range("B1").Value = "Title"
$ws2.Range("D1").Copy($ws3.Range("C1"))
$ws3.Range("C1").Value = "Deleteaction"

$ws3.Range("B2").Value = "Selenium"
$ws3.Range("C2").Value = "Yes"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "http://192.168.99.100:5000/Home/ServiceListing/?id=60d82811b4b34c00014c766a")
$ws3.Range("A2").Style = "Hyperlink"

$ws3.Columns.Item(1).ColumnWidth = 68.3
$ws3.Columns.Item(2).ColumnWidth = 16.5
$ws3.Columns.Item(3).ColumnWidth = 27

$ws3.Range("D8").Select()

# -----------------------------------------------------------------
# New sheet: ServiceDetail
# -----------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "ServiceDetail"

$ws4.Range("A1").Value = "Chatboxvalues"
$ws4.Range("B1").Value = "url"
$ws4.Range("A2").Value = "Hi How are you?"
$ws4.Range("B2").Value = "Selenium"

$ws4.Hyperlinks.Add($ws4.Range("B2"), "http://192.168.99.100:5000/Home/ServiceDetail?id=60d460539b4eae0001de9f70")
$ws4.Range("B2").Style = "Hyperlink"

$ws4.Columns.Item(1).ColumnWidth = 23.6
$ws4.Columns.Item(2).ColumnWidth = 69.6

$ws4.Range("B20").Select()

# -----------------------------------------------------------------
# Leave ShareSkill as the active/visible tab, matching the target
# -----------------------------------------------------------------
$ws2.Activate()
$ws2.Range("P10").Select()
